$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
# Row 40
$ws.Range("H40").Value = 2859.4
$ws.Range("J40").Value = 2649.5
$ws.Range("L40").Value = 2649.5
$ws.Range("N40").Value = -2999.5
# Row 76
$ws.Range("H76").Value = 43487256
$ws.Range("I76").Value = 66676468
$ws.Range("J76").Value = 7480.875
$ws.Range("K76").Value = 66676468
$ws.Range("L76").Value = 7480.875
$ws.Range("M76").Value = -66676153
$ws.Range("N76").Value = -8110.875
# Row 79
$ws.Range("H79").Value = 43487256
$ws.Range("I79").Value = 66676468
$ws.Range("J79").Value = 7480.875
$ws.Range("K79").Value = 66676468
$ws.Range("L79").Value = 7480.875
$ws.Range("M79").Value = -66675376
$ws.Range("N79").Value = -9664.875
# Row 106
$ws.Range("H106").Value = 1490.1875
$ws.Range("I106").Value = 1491.0333
$ws.Range("J106").Value = 1477.5
$ws.Range("K106").Value = 1491.0333
$ws.Range("L106").Value = 1477.5
$ws.Range("M106").Value = -860.0333000000001
$ws.Range("N106").Value = -2739.5
# Row 137
$ws.Range("H137").Value = 6533.5806
$ws.Range("I137").Value = 3357.5293
$ws.Range("J137").Value = 10390.214
$ws.Range("K137").Value = 10072.5879
$ws.Range("L137").Value = 31170.642
$ws.Range("M137").Value = -7522.5879
$ws.Range("N137").Value = -36270.642

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 1356594.8
$ws.Range("I32").Value = 795.44446
$ws.Range("J32").Value = 11118350
$ws.Range("K32").Value = 795.44446
$ws.Range("L32").Value = 11118350
$ws.Range("M32").Value = -508.44446
$ws.Range("N32").Value = -11118924
# Row 132
$ws.Range("H132").Value = 847669
$ws.Range("I132").Value = 1062178
$ws.Range("K132").Value = 3186534
$ws.Range("M132").Value = -3184004

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
# Row 86
$ws.Range("H86").Value = 6324.1377
$ws.Range("I86").Value = 5170.2104
$ws.Range("J86").Value = 8516.6
$ws.Range("K86").Value = 5170.2104
$ws.Range("L86").Value = 8516.6
$ws.Range("M86").Value = -4047.2104
$ws.Range("N86").Value = -10762.6
# Row 89
$ws.Range("H89").Value = 6324.1377
$ws.Range("I89").Value = 5170.2104
$ws.Range("J89").Value = 8516.6
$ws.Range("K89").Value = 25851.052
$ws.Range("L89").Value = 42583
$ws.Range("M89").Value = -20235.052
$ws.Range("N89").Value = -53815

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 7638.5
$ws.Range("I31").Value = 3099.2
$ws.Range("J31").Value = 9701.817999999999
$ws.Range("K31").Value = 3099.2
$ws.Range("L31").Value = 9701.817999999999
$ws.Range("M31").Value = -2804.2
$ws.Range("N31").Value = -10291.818
# Row 34
$ws.Range("H34").Value = 7638.5
$ws.Range("I34").Value = 3099.2
$ws.Range("J34").Value = 9701.817999999999
$ws.Range("K34").Value = 3099.2
$ws.Range("L34").Value = 9701.817999999999
$ws.Range("M34").Value = -2897.2
$ws.Range("N34").Value = -10105.818
# Row 99
$ws.Range("H99").Value = 8267325
$ws.Range("I99").Value = 12988768
$ws.Range("J99").Value = 4799.75
$ws.Range("K99").Value = 12988768
$ws.Range("L99").Value = 4799.75
$ws.Range("M99").Value = -12987270
$ws.Range("N99").Value = -7795.75
# Row 126
$ws.Range("H126").Value = 8267325
$ws.Range("I126").Value = 12988768
$ws.Range("J126").Value = 4799.75
$ws.Range("K126").Value = 38966304
$ws.Range("L126").Value = 14399.25
$ws.Range("M126").Value = -38963834
$ws.Range("N126").Value = -19339.25
# Row 134
$ws.Range("H134").Value = 37044904
$ws.Range("I134").Value = 166677250
$ws.Range("J134").Value = 7092.8096
$ws.Range("K134").Value = 500031750
$ws.Range("L134").Value = 21278.4288
$ws.Range("M134").Value = -500029215
$ws.Range("N134").Value = -26348.4288

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
# Row 34
$ws.Range("H34").Value = 2003.2727
$ws.Range("I34").Value = 237.4
$ws.Range("J34").Value = 3474.8333
$ws.Range("K34").Value = 712.2
$ws.Range("L34").Value = 10424.4999
$ws.Range("M34").Value = -628.2
$ws.Range("N34").Value = -10592.4999
# Row 39
$ws.Range("H39").Value = 6162.5
$ws.Range("J39").Value = 6162.5
$ws.Range("L39").Value = 18487.5
$ws.Range("N39").Value = -19075.5
# Row 68
$ws.Range("H68").Value = 174689.83
$ws.Range("I68").Value = 1333.2
$ws.Range("J68").Value = 360429.06
$ws.Range("K68").Value = 3999.6
$ws.Range("L68").Value = 1081287.18
$ws.Range("M68").Value = -3188.6
$ws.Range("N68").Value = -1082909.18
# Row 71
$ws.Range("H71").Value = 174689.83
$ws.Range("I71").Value = 1333.2
$ws.Range("J71").Value = 360429.06
$ws.Range("K71").Value = 11998.8
$ws.Range("L71").Value = 3243861.54
$ws.Range("M71").Value = -7942.800000000001
$ws.Range("N71").Value = -3251973.54
# Row 114
$ws.Range("H114").Value = 316.16666
$ws.Range("I114").Value = 316.16666
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 948.4999799999999
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 2305.50002
$ws.Range("N114").ClearContents()
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
# Row 139
$ws.Range("H139").Value = 115312.375
$ws.Range("I139").Value = 3214.1428
$ws.Range("J139").Value = 900000
$ws.Range("K139").Value = 9642.428400000001
$ws.Range("L139").Value = 2700000
$ws.Range("M139").Value = -4502.428400000001
$ws.Range("N139").Value = -2710280

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
# Row 102
$ws.Range("H102").Value = 4395.5
$ws.Range("I102").Value = 2798.5144
$ws.Range("J102").Value = 8121.8
$ws.Range("K102").Value = 2798.5144
$ws.Range("L102").Value = 8121.8
$ws.Range("M102").Value = -1176.5144
$ws.Range("N102").Value = -11365.8

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
# Row 68
$ws.Range("H68").Value = 1560.4
$ws.Range("I68").Value = 1500.5
$ws.Range("J68").Value = 1800
$ws.Range("K68").Value = 1500.5
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = -751.5
$ws.Range("N68").Value = -3298
# Row 71
$ws.Range("H71").Value = 1560.4
$ws.Range("I71").Value = 1500.5
$ws.Range("J71").Value = 1800
$ws.Range("K71").Value = 7502.5
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -3758.5
$ws.Range("N71").Value = -16488
# Row 93
$ws.Range("H93").Value = 1549.5454
$ws.Range("I93").Value = 4045
$ws.Range("K93").Value = 4045
$ws.Range("M93").Value = -2797
# Row 132
$ws.Range("H132").Value = 6478.8
$ws.Range("I132").Value = 6192.6523
$ws.Range("K132").Value = 18577.9569
$ws.Range("M132").Value = -16047.9569

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
# Row 107
$ws.Range("H107").Value = 1265.3704
$ws.Range("I107").Value = 1430.6364
$ws.Range("K107").Value = 4291.9092
$ws.Range("M107").Value = -2371.9092
# Row 136
$ws.Range("H136").Value = 8339966.5
$ws.Range("I136").Value = 10208310
$ws.Range("J136").Value = 17345.455
$ws.Range("K136").Value = 30624930
$ws.Range("L136").Value = 52036.36500000001
$ws.Range("M136").Value = -30622380
$ws.Range("N136").Value = -57136.36500000001
